# Change dashes to underscores in preparation for shinyfilter
# (applies to the Filename values in column A)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    if ($current -ne $null) {
        $updated = $current.Replace("-", "_")
        if ($updated -ne $current) {
            $cell.Value = $updated
        }
    }
}
